# Toevoegen CORS aan logboek
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (cloning row 10's formatting, which already
# has the A/B/C date-time layout used by other logbook entries), pushing the
# existing rows 11-14 down to 12-15.
$ws.Rows("10:10").Copy()
$ws.Rows("11:11").Insert()

# Fill in the new logbook entry in the freshly inserted row 11. B/C keep the
# row's existing date/time number formats, so write the raw serial values
# (46018 = 2025-12-27, 0.9375 = 22:30) instead of DateTime objects, which
# would otherwise stamp a brand new number format onto the cell.
$ws.Range("A11").Value = "CORS"
$ws.Range("B11").Value = 46018
$ws.Range("C11").Value = 0.9375

# Update the active selection to mirror the authored state (D13).
$ws.Range("D13").Select()
